$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C..G to D..H)
$ws.Columns("C:C").Insert()

# Column widths for the new/changed columns
$ws.Columns("B:B").ColumnWidth = 24.1
$ws.Columns("C:C").ColumnWidth = 17.59

# New row with the requirement link (set first so shared-string order matches)
$ws.Range("B5").Value = "Del 2 – Inloggningsfunktion"

# New header cell in the inserted column
$ws.Range("C4").Value = "Req link"

$ws.Range("C5").Value = "Req Link"
$ws.Hyperlinks.Add($ws.Range("C5"), "https://example.com")

# Update the selection to match the new active cell
$ws.Range("C5").Select()
